# DOC(X) filter: BAF -> BAILS in the smart tag mapping
#
# Rewrites the three body paragraphs so that:
#  - each paragraph gets an explicit pPr/pStyle="Normal"
#  - the ordinal word ("1st"/"2nd"/"3rd") is split into its own runs, with
#    the ordinal suffix raised to superscript
#  - the smart tag's RDF attribute names are remapped from the old
#    TSCP "BAF" URNs to the BAILS URNs (and a StartValidity attribute
#    replaces the old BusinessAuthorizationDate one)
#  - the stray "_GoBack" bookmark pair is dropped

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$p1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr/>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>st</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, non-business.</w:t>
  </w:r>
</w:p>
'@
$p1.Range.InsertXML($p1Xml)

$p2 = $d.Paragraphs(2)
$p2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:smartTag w:uri="http://www.w3.org/1999/02/22-rdf-syntax-ns#" w:element="RDF">
    <w:smartTagPr>
      <w:attr w:name="urn:bails:ExportControl:Authorization:StartValidity" w:val="2015-11-27"/>
      <w:attr w:name="urn:bails:ExportControl:BusinessAuthorization:Identifier" w:val="urn:example:tscp:1"/>
      <w:attr w:name="urn:bails:ExportControl:BusinessAuthorizationCategory:Identifier" w:val="urn:example:tscp:1:confidential"/>
    </w:smartTagPr>
  </w:smartTag>
  <w:r>
    <w:rPr/>
    <w:t>2</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>nd</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, confidential.</w:t>
  </w:r>
</w:p>
'@
$p2.Range.InsertXML($p2Xml)

$p3 = $d.Paragraphs(3)
$p3Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr/>
    <w:t>3</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>rd</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, non-business.</w:t>
  </w:r>
</w:p>
'@
$p3.Range.InsertXML($p3Xml)
